# This workbook's data rows (2,3,5,6,7) got cyclically permuted during the
# source update: the record that used to live in row 2 now lives in row 5,
# row 5's old record moved to row 3, row 3's old record moved to row 7,
# row 7's old record moved to row 6, and row 6's old record moved to row 2.
# (Row 4 is untouched.) Apply the resulting cell values directly per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (now holds what used to be row 6's record) ---
$ws.Range("A2").Value = 111402344
$ws.Range("B2").Value = 90666
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 4364
$ws.Range("F2").Value = "Dropptaggsvamp"
$ws.Range("G2").Value = "Hydnellum ferrugineum"
$ws.Range("H2").Value = "(Fr.:Fr.) P. Karst."
# This record has no Ålder-Stadium / Kön / Aktivitet / Metod values, so clear them.
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("Q2").Value = 545197.7500373307
$ws.Range("R2").Value = 7020179.372318991

# --- Row 3 (now holds what used to be row 5's record) ---
$ws.Range("A3").Value = 111402339
$ws.Range("B3").Value = 77267
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 545198.1129081531
$ws.Range("R3").Value = 7020057.514641967

# --- Row 5 (now holds what used to be row 2's record) ---
$ws.Range("A5").Value = 111402340
$ws.Range("B5").Value = 56414
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
# This record carries Ålder-Stadium / Kön / Aktivitet / Metod info.
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "spel/sång"
$ws.Range("N5").Value = ""

# --- Row 6 (now holds what used to be row 7's record) ---
$ws.Range("A6").Value = 111402337
$ws.Range("B6").Value = 96370
$ws.Range("E6").Value = 219847
$ws.Range("F6").Value = "Tvåblad"
$ws.Range("G6").Value = "Neottia ovata"
$ws.Range("H6").Value = "(L.) Buff. & Fingerh."
$ws.Range("Q6").Value = 545198.1129081531
$ws.Range("R6").Value = 7020057.514641967

# --- Row 7 (now holds what used to be row 3's record) ---
$ws.Range("A7").Value = 111402343
$ws.Range("B7").Value = 77186
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 353
$ws.Range("F7").Value = "Dvärgbägarlav"
$ws.Range("G7").Value = "Cladonia parasitica"
$ws.Range("H7").Value = "(Hoffm.) Hoffm."
$ws.Range("Q7").Value = 545197.7500373307
$ws.Range("R7").Value = 7020179.372318991
